# Update ticket_id values (column D) and corresponding ticket_link URLs (column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseUrl = "https://8000-gaius1-qrcodeticketings-256zdob18db.ws-eu63.gitpod.io/ticket/"

$updates = @{
    2 = 582050
    3 = 578600
    4 = 322289
    5 = 894139
    6 = 718113
    7 = 185950
}

foreach ($row in $updates.Keys) {
    $newId = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $newId
    $ws.Cells.Item($row, 5).Value = "$baseUrl$newId"
}
